$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1964017991004498
$ws.Range("C2").Value = 0.5547226386806596
$ws.Range("J2").Value = 0.01799100449775112
$ws.Range("P2").Value = 0.1334332833583208
$ws.Range("S2").Value = 0.09745127436281859
$ws.Range("B3").Value = 0.01038961038961039
$ws.Range("C3").Value = 0.02337662337662338
$ws.Range("J3").Value = 0.02597402597402598
$ws.Range("P3").Value = 0.7454545454545455
$ws.Range("S3").Value = 0.1948051948051948
$ws.Range("J4").Value = 0.03703703703703703
$ws.Range("P4").Value = 0.6790123456790124
$ws.Range("S4").Value = 0.2839506172839506
$ws.Range("B6").Value = 0.06791569086651054
$ws.Range("D6").Value = 0.00702576112412178
$ws.Range("F6").Value = 0.05152224824355972
$ws.Range("J6").Value = 0.2646370023419204
$ws.Range("O6").Value = 0.01873536299765808
$ws.Range("Q6").Value = 0.1639344262295082
$ws.Range("R6").Value = 0.07494145199063232
$ws.Range("S6").Value = 0.351288056206089
$ws.Range("B7").Value = 0.0970873786407767
$ws.Range("D7").Value = 0.02184466019417476
$ws.Range("E7").Value = 0.002427184466019417
$ws.Range("F7").Value = 0.0412621359223301
$ws.Range("J7").Value = 0.1747572815533981
$ws.Range("O7").Value = 0.01941747572815534
$ws.Range("Q7").Value = 0.1699029126213592
$ws.Range("R7").Value = 0.06796116504854369
$ws.Range("S7").Value = 0.4053398058252427
$ws.Range("B8").Value = 0.118510158013544
$ws.Range("D8").Value = 0.02031602708803612
$ws.Range("F8").Value = 0.08916478555304741
$ws.Range("J8").Value = 0.1060948081264108
$ws.Range("O8").Value = 0.01467268623024831
$ws.Range("Q8").Value = 0.1783295711060948
$ws.Range("R8").Value = 0.08803611738148984
$ws.Range("S8").Value = 0.3848758465011287
$ws.Range("B9").Value = 0.08950617283950617
$ws.Range("D9").Value = 0.006172839506172839
$ws.Range("E9").Value = 0.00308641975308642
$ws.Range("F9").Value = 0.05864197530864197
$ws.Range("J9").Value = 0.1327160493827161
$ws.Range("O9").Value = 0.01851851851851852
$ws.Range("Q9").Value = 0.191358024691358
$ws.Range("R9").Value = 0.06790123456790123
$ws.Range("S9").Value = 0.4320987654320987
$ws.Range("B10").Value = 0.1176470588235294
$ws.Range("D10").Value = 0.01863354037267081
$ws.Range("E10").Value = 0.001461454146876142
$ws.Range("F10").Value = 0.06393861892583121
$ws.Range("J10").Value = 0.1293386919985385
$ws.Range("O10").Value = 0.01607599561563756
$ws.Range("Q10").Value = 0.2246985750822068
$ws.Range("R10").Value = 0.08111070515162587
$ws.Range("S10").Value = 0.3470953598830837
$ws.Range("G11").Value = 0.1526479750778816
$ws.Range("J11").Value = 0.08411214953271028
$ws.Range("K11").Value = 0.1962616822429906
$ws.Range("L11").Value = 0.5607476635514018
$ws.Range("S11").Value = 0.006230529595015576
$ws.Range("G12").Value = 0.7142857142857143
$ws.Range("J12").Value = 0.2371967654986523
$ws.Range("K12").Value = 0.002695417789757413
$ws.Range("L12").Value = 0.02425876010781671
$ws.Range("S12").Value = 0.0215633423180593
$ws.Range("F13").Value = 0.01098901098901099
$ws.Range("G13").Value = 0.6373626373626373
$ws.Range("J13").Value = 0.2857142857142857
$ws.Range("S13").Value = 0.06593406593406594
$ws.Range("F15").Value = 0.0101010101010101
$ws.Range("H15").Value = 0.1338383838383838
$ws.Range("I15").Value = 0.06818181818181818
$ws.Range("J15").Value = 0.3762626262626262
$ws.Range("K15").Value = 0.06565656565656566
$ws.Range("M15").Value = 0.02272727272727273
$ws.Range("O15").Value = 0.04797979797979798
$ws.Range("S15").Value = 0.2752525252525252
$ws.Range("F16").Value = 0.01431980906921241
$ws.Range("H16").Value = 0.1599045346062052
$ws.Range("I16").Value = 0.06921241050119331
$ws.Range("J16").Value = 0.4391408114558473
$ws.Range("K16").Value = 0.09785202863961814
$ws.Range("M16").Value = 0.01909307875894988
$ws.Range("N16").Value = 0.002386634844868735
$ws.Range("O16").Value = 0.03341288782816229
$ws.Range("S16").Value = 0.1646778042959427
$ws.Range("F17").Value = 0.01131687242798354
$ws.Range("H17").Value = 0.1851851851851852
$ws.Range("I17").Value = 0.0668724279835391
$ws.Range("J17").Value = 0.4403292181069959
$ws.Range("K17").Value = 0.09670781893004116
$ws.Range("M17").Value = 0.01851851851851852
$ws.Range("O17").Value = 0.04938271604938271
$ws.Range("S17").Value = 0.1316872427983539
$ws.Range("F18").Value = 0.01058201058201058
$ws.Range("H18").Value = 0.1851851851851852
$ws.Range("I18").Value = 0.09259259259259259
$ws.Range("J18").Value = 0.3941798941798942
$ws.Range("K18").Value = 0.1058201058201058
$ws.Range("M18").Value = 0.02116402116402116
$ws.Range("O18").Value = 0.08201058201058201
$ws.Range("S18").Value = 0.1084656084656085
$ws.Range("F19").Value = 0.01498422712933754
$ws.Range("H19").Value = 0.2058359621451104
$ws.Range("I19").Value = 0.06703470031545741
$ws.Range("J19").Value = 0.3907728706624606
$ws.Range("K19").Value = 0.1214511041009464
$ws.Range("M19").Value = 0.02129337539432177
$ws.Range("N19").Value = 0.0003943217665615142
$ws.Range("O19").Value = 0.05954258675078864
$ws.Range("S19").Value = 0.1186908517350158
